$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2144638403990025
$ws.Cells.Item(2, 3).Value = 0.5236907730673317
$ws.Cells.Item(2, 10).Value = 0.01745635910224439
$ws.Cells.Item(2, 16).Value = 0.1596009975062344
$ws.Cells.Item(2, 19).Value = 0.08478802992518704
$ws.Cells.Item(3, 2).Value = 0.009259259259259259
$ws.Cells.Item(3, 3).Value = 0.01851851851851852
$ws.Cells.Item(3, 10).Value = 0.02314814814814815
$ws.Cells.Item(3, 16).Value = 0.7083333333333334
$ws.Cells.Item(3, 19).Value = 0.2407407407407407
$ws.Cells.Item(4, 10).Value = 0.08333333333333333
$ws.Cells.Item(4, 16).Value = 0.7291666666666666
$ws.Cells.Item(4, 19).Value = 0.1875
$ws.Cells.Item(6, 2).Value = 0.0546218487394958
$ws.Cells.Item(6, 4).Value = 0.008403361344537815
$ws.Cells.Item(6, 6).Value = 0.03781512605042017
$ws.Cells.Item(6, 10).Value = 0.3067226890756303
$ws.Cells.Item(6, 17).Value = 0.1134453781512605
$ws.Cells.Item(6, 18).Value = 0.09243697478991597
$ws.Cells.Item(6, 19).Value = 0.3865546218487395
$ws.Cells.Item(7, 2).Value = 0.1372549019607843
$ws.Cells.Item(7, 4).Value = 0.02352941176470588
$ws.Cells.Item(7, 6).Value = 0.04705882352941176
$ws.Cells.Item(7, 10).Value = 0.09019607843137255
$ws.Cells.Item(7, 15).Value = 0.0196078431372549
$ws.Cells.Item(7, 17).Value = 0.1882352941176471
$ws.Cells.Item(7, 18).Value = 0.08235294117647059
$ws.Cells.Item(7, 19).Value = 0.4117647058823529
$ws.Cells.Item(8, 2).Value = 0.1064814814814815
$ws.Cells.Item(8, 4).Value = 0.01851851851851852
$ws.Cells.Item(8, 5).Value = 0.002314814814814815
$ws.Cells.Item(8, 6).Value = 0.05787037037037037
$ws.Cells.Item(8, 10).Value = 0.1273148148148148
$ws.Cells.Item(8, 15).Value = 0.02546296296296296
$ws.Cells.Item(8, 17).Value = 0.1875
$ws.Cells.Item(8, 18).Value = 0.09490740740740741
$ws.Cells.Item(8, 19).Value = 0.3796296296296297
$ws.Cells.Item(9, 2).Value = 0.0752212389380531
$ws.Cells.Item(9, 4).Value = 0.02212389380530973
$ws.Cells.Item(9, 6).Value = 0.05752212389380531
$ws.Cells.Item(9, 10).Value = 0.1238938053097345
$ws.Cells.Item(9, 15).Value = 0.03539823008849557
$ws.Cells.Item(9, 17).Value = 0.2123893805309734
$ws.Cells.Item(9, 18).Value = 0.08849557522123894
$ws.Cells.Item(9, 19).Value = 0.3849557522123894
$ws.Cells.Item(10, 2).Value = 0.1373937677053824
$ws.Cells.Item(10, 4).Value = 0.01912181303116147
$ws.Cells.Item(10, 6).Value = 0.06303116147308782
$ws.Cells.Item(10, 10).Value = 0.1232294617563739
$ws.Cells.Item(10, 15).Value = 0.02407932011331445
$ws.Cells.Item(10, 17).Value = 0.231586402266289
$ws.Cells.Item(10, 18).Value = 0.0594900849858357
$ws.Cells.Item(10, 19).Value = 0.3420679886685553
$ws.Cells.Item(11, 7).Value = 0.1679790026246719
$ws.Cells.Item(11, 10).Value = 0.06824146981627296
$ws.Cells.Item(11, 11).Value = 0.1837270341207349
$ws.Cells.Item(11, 12).Value = 0.5590551181102362
$ws.Cells.Item(11, 19).Value = 0.02099737532808399
$ws.Cells.Item(12, 7).Value = 0.7880184331797235
$ws.Cells.Item(12, 10).Value = 0.1612903225806452
$ws.Cells.Item(12, 11).Value = 0.01382488479262673
$ws.Cells.Item(12, 12).Value = 0.0184331797235023
$ws.Cells.Item(12, 19).Value = 0.0184331797235023
$ws.Cells.Item(13, 7).Value = 0.574468085106383
$ws.Cells.Item(13, 10).Value = 0.3191489361702128
$ws.Cells.Item(13, 19).Value = 0.1063829787234043
$ws.Cells.Item(15, 6).Value = 0.02145922746781116
$ws.Cells.Item(15, 8).Value = 0.1502145922746781
$ws.Cells.Item(15, 9).Value = 0.07725321888412018
$ws.Cells.Item(15, 10).Value = 0.3218884120171674
$ws.Cells.Item(15, 11).Value = 0.06866952789699571
$ws.Cells.Item(15, 13).Value = 0.004291845493562232
$ws.Cells.Item(15, 14).Value = 0.004291845493562232
$ws.Cells.Item(15, 15).Value = 0.04291845493562232
$ws.Cells.Item(15, 19).Value = 0.3090128755364807
$ws.Cells.Item(16, 6).Value = 0.04471544715447155
$ws.Cells.Item(16, 8).Value = 0.1382113821138211
$ws.Cells.Item(16, 9).Value = 0.07317073170731707
$ws.Cells.Item(16, 10).Value = 0.4471544715447154
$ws.Cells.Item(16, 11).Value = 0.1300813008130081
$ws.Cells.Item(16, 13).Value = 0.01219512195121951
$ws.Cells.Item(16, 15).Value = 0.04878048780487805
$ws.Cells.Item(16, 19).Value = 0.1056910569105691
$ws.Cells.Item(17, 6).Value = 0.02646502835538752
$ws.Cells.Item(17, 8).Value = 0.1398865784499055
$ws.Cells.Item(17, 9).Value = 0.1209829867674858
$ws.Cells.Item(17, 10).Value = 0.4177693761814745
$ws.Cells.Item(17, 11).Value = 0.1077504725897921
$ws.Cells.Item(17, 13).Value = 0.0113421550094518
$ws.Cells.Item(17, 15).Value = 0.0661625708884688
$ws.Cells.Item(17, 19).Value = 0.109640831758034
$ws.Cells.Item(18, 6).Value = 0.03225806451612903
$ws.Cells.Item(18, 8).Value = 0.1397849462365591
$ws.Cells.Item(18, 9).Value = 0.05376344086021505
$ws.Cells.Item(18, 10).Value = 0.4516129032258064
$ws.Cells.Item(18, 11).Value = 0.1129032258064516
$ws.Cells.Item(18, 13).Value = 0.02150537634408602
$ws.Cells.Item(18, 15).Value = 0.08064516129032258
$ws.Cells.Item(18, 19).Value = 0.1075268817204301
$ws.Cells.Item(19, 6).Value = 0.01827485380116959
$ws.Cells.Item(19, 8).Value = 0.1988304093567251
$ws.Cells.Item(19, 9).Value = 0.08406432748538012
$ws.Cells.Item(19, 10).Value = 0.3625730994152047
$ws.Cells.Item(19, 11).Value = 0.125
$ws.Cells.Item(19, 13).Value = 0.02631578947368421
$ws.Cells.Item(19, 14).Value = 0.0007309941520467836
$ws.Cells.Item(19, 15).Value = 0.0577485380116959
$ws.Cells.Item(19, 19).Value = 0.1264619883040936